$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.027.59'
$ws.Range("E2").Value = '  -0.55%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.428.59'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.20%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '409.61'
$ws.Range("E5").Value = '  +0.72%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '129.06'
$ws.Range("E6").Value = '  -2.83%  '
$ws.Range("E7").Value = '  +6.33%  '
$ws.Range("E8").Value = '  -0.10%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.741'
$ws.Range("E9").Value = '  +7.54%  '
$ws.Range("E10").Value = '  +5.52%  '
$ws.Range("E11").Value = '  +2.08%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000224'
$ws.Range("E12").Value = '  +50.83%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '9.21'
$ws.Range("E13").Value = '  +9.58%  '
$ws.Range("E15").Value = '  +7.96%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.969.04'
$ws.Range("E16").Value = '  -0.06%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.444.53'
$ws.Range("E17").Value = '  +0.39%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.59'
$ws.Range("E18").Value = '  +7.46%  '
$ws.Range("E19").Value = '  +7.58%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '62.010.25'
$ws.Range("E20").Value = '  -0.44%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '452.41'
$ws.Range("E21").Value = '  +45.12%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '91.76'
$ws.Range("E22").Value = '  +9.09%  '
$ws.Range("E23").Value = '  +1.43%  '
$ws.Range("E24").Value = '  +2.05%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.25'
$ws.Range("E25").Value = '  +2.26%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '32.98'
$ws.Range("E26").Value = '  +10.85%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.90'
$ws.Range("E27").Value = '  +9.31%  '
$ws.Range("E28").Value = '  +0.77%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.70'
$ws.Range("E29").Value = '  -1.64%  '
$ws.Range("E30").Value = '  -0.65%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '12.00'
$ws.Range("E31").Value = '  +5.82%  '
$ws.Range("E32").Value = '  -0.96%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '43.16'
$ws.Range("E33").Value = '  -1.60%  '
$ws.Range("E34").Value = '  -0.45%  '
$ws.Range("E35").Value = '  -0.08%  '
$ws.Range("E36").Value = '  +3.09%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '54.37'
$ws.Range("E37").Value = '  +5.30%  '
$ws.Range("E38").Value = '  -0.15%  '
$ws.Range("E39").Value = '  +1.36%  '
$ws.Range("E40").Value = '  +7.68%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.322'
$ws.Range("E41").Value = '  +1.34%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.94'
$ws.Range("E42").Value = '  -2.05%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '142.57'
$ws.Range("E43").Value = '  +0.20%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.27'
$ws.Range("E44").Value = '  +8.93%  '
$ws.Range("E45").Value = '  +1.02%  '
$ws.Range("E46").Value = '  +13.33%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '16.67'
$ws.Range("E47").Value = '  -0.69%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '22.34'
$ws.Range("E48").Value = '  +5.02%  '
$ws.Range("E49").Value = '  +9.40%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.774.87'
$ws.Range("E50").Value = '  +0.06%  '
$ws.Range("E51").Value = '  +15.49%  '
